$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing "Branch" column (J, rows 1-27 -- the instruction table)
# one column to the right, into K, preserving its values and highlight
# styling. Row 28 (the legend/footnote row) is left untouched.
$ws.Range("J1:J27").Copy()
$ws.Range("K1:K27").PasteSpecial(-4104)
$ws.Application.CutCopyMode = $false

# J now becomes the new "JPR" column. Header first.
$ws.Range("J1").Value = "JPR"

# Every instruction row defaults to 0 (plain, unstyled) in the new column.
foreach ($r in 2..27) {
    $cell = $ws.Range("J$r")
    $cell.ClearFormats()
    $cell.Value = 0
}

# JPR and JRL (rows 23 and 24) are flagged in the new JPR column instead of
# the generic JP column: JP (I) goes back to plain 0, and the new JPR cell
# (J) becomes a highlighted 1 -- the same "true flag" look used elsewhere
# on the sheet (light-yellow fill + thin grey border).
foreach ($r in 23, 24) {
    $jp = $ws.Range("I$r")
    $jp.ClearFormats()
    $jp.Value = 0

    $jpr = $ws.Range("J$r")
    $jpr.Value = 1
    $jpr.Interior.Color = 13434879
    $jpr.Interior.Pattern = 1
    $jpr.Borders.LineStyle = 1
    $jpr.Borders.Color = 11711154
}

# Match the edited sheet's selection.
$ws.Range("A15:G15").Select()
